$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.945.80'
Set-TextValue $ws.Range('E2') '  +1.01%  '
Set-TextValue $ws.Range('D3') '1.636.76'
Set-TextValue $ws.Range('E3') '  +1.81%  '
Set-TextValue $ws.Range('E4') '  +0.39%  '
Set-TextValue $ws.Range('D5') '214.99'
Set-TextValue $ws.Range('E5') '  +1.29%  '
Set-TextValue $ws.Range('D6') '0.519'
Set-TextValue $ws.Range('E6') '  +0.25%  '
Set-TextValue $ws.Range('D7') '0.999'
Set-TextValue $ws.Range('E7') '  +0.38%  '
Set-TextValue $ws.Range('D8') '28.85'
Set-TextValue $ws.Range('E8') '  -0.38%  '
Set-TextValue $ws.Range('E9') '  +0.00%  '
Set-TextValue $ws.Range('E10') '  +0.33%  '
Set-TextValue $ws.Range('E11') '  -0.30%  '
Set-TextValue $ws.Range('D12') '1.869.40'
Set-TextValue $ws.Range('E12') '  +1.69%  '
Set-TextValue $ws.Range('D13') '1.638.68'
Set-TextValue $ws.Range('E13') '  +2.13%  '
Set-TextValue $ws.Range('E14') '  -0.18%  '
Set-TextValue $ws.Range('D15') '9.33'
Set-TextValue $ws.Range('E15') '  +11.60%  '
Set-TextValue $ws.Range('D16') '29.960.99'
Set-TextValue $ws.Range('E16') '  +0.99%  '
Set-TextValue $ws.Range('D17') '3.86'
Set-TextValue $ws.Range('E17') '  +1.11%  '
Set-TextValue $ws.Range('D18') '64.27'
Set-TextValue $ws.Range('E18') '  -0.46%  '
Set-TextValue $ws.Range('D19') '242.04'
Set-TextValue $ws.Range('E19') '  +0.24%  '
Set-TextValue $ws.Range('D20') '0.0₃0703'
Set-TextValue $ws.Range('E20') '  +0.05%  '
Set-TextValue $ws.Range('E21') '  +0.32%  '
Set-TextValue $ws.Range('D22') '4.14'
Set-TextValue $ws.Range('E22') '  +2.20%  '
Set-TextValue $ws.Range('D23') '9.81'
Set-TextValue $ws.Range('E23') '  +2.97%  '
Set-TextValue $ws.Range('D24') '2.17'
Set-TextValue $ws.Range('E24') '  +3.14%  '
Set-TextValue $ws.Range('D25') '158.09'
Set-TextValue $ws.Range('E25') '  +0.96%  '
Set-TextValue $ws.Range('E26') '  -0.41%  '
Set-TextValue $ws.Range('E27') '  +0.62%  '
Set-TextValue $ws.Range('D28') '6.59'
Set-TextValue $ws.Range('E28') '  +0.93%  '
Set-TextValue $ws.Range('E29') '  +0.33%  '
Set-TextValue $ws.Range('E30') '  +2.08%  '
Set-TextValue $ws.Range('E31') '  +3.81%  '
Set-TextValue $ws.Range('E32') '  +3.44%  '
Set-TextValue $ws.Range('E33') '  +0.27%  '
Set-TextValue $ws.Range('D34') '1.431.38'
Set-TextValue $ws.Range('E34') '  +0.34%  '
Set-TextValue $ws.Range('D35') '1.67'
Set-TextValue $ws.Range('E35') '  +5.17%  '
Set-TextValue $ws.Range('E36') '  -1.80%  '
Set-TextValue $ws.Range('D37') '2.79'
Set-TextValue $ws.Range('E37') '  -2.70%  '
Set-TextValue $ws.Range('B38') 'HuobiToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D38') '2.29'
Set-TextValue $ws.Range('E38') '  +0.36%  '
Set-TextValue $ws.Range('B39') 'VeChain'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D39') '0.0171'
Set-TextValue $ws.Range('E39') '  +0.90%  '
Set-TextValue $ws.Range('D40') '75.84'
Set-TextValue $ws.Range('E40') '  +10.85%  '
Set-TextValue $ws.Range('D41') '0.554'
Set-TextValue $ws.Range('E41') '  -0.02%  '
Set-TextValue $ws.Range('B42') 'Kaspa'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D42') '0.0504'
Set-TextValue $ws.Range('E42') '  +1.42%  '
Set-TextValue $ws.Range('B43') 'RenderToken'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D43') '1.99'
Set-TextValue $ws.Range('E43') '  +0.96%  '
Set-TextValue $ws.Range('B44') 'ARBITRUM'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D44') '0.830'
Set-TextValue $ws.Range('E44') '  +0.74%  '
Set-TextValue $ws.Range('E45') '  +0.40%  '
Set-TextValue $ws.Range('E46') '  +0.84%  '
Set-TextValue $ws.Range('D47') '51.26'
Set-TextValue $ws.Range('E47') '  -5.86%  '
Set-TextValue $ws.Range('B48') 'RocketPoolETH'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D48') '1.776.48'
Set-TextValue $ws.Range('E48') '  +1.67%  '
Set-TextValue $ws.Range('B49') 'FraxShare'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D49') '5.35'
Set-TextValue $ws.Range('E49') '  -1.25%  '
Set-TextValue $ws.Range('B50') 'Quant'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D50') '90.55'
Set-TextValue $ws.Range('E50') '  +3.98%  '
Set-TextValue $ws.Range('B51') 'BabyDogeCoin'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D51') '0.0₆0113'
Set-TextValue $ws.Range('E51') '  +10.86%  '
